# Update column F (dSF) values on Sheet1 to match the repulled / pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    6  = -4
    8  = -2
    11 = -4
    15 = -5
    18 = 8
    19 = -5
    21 = -11
    24 = 0
    27 = -3
    29 = 0
    31 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
